# This workbook contains one data table (rows 2-35) describing daily price
# observations for Alcachofa (artichoke). The edit re-shuffles the
# per-observation columns (Fecha, Variedad, Calidad, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Origen, Precio $/Kg) across the
# existing 34 rows - i.e. row N now shows the values that used to belong to
# a different row (a pure permutation of those columns; the rest of the
# columns - Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría,
# Unidad de comercialización, Kg o Unidades, Clasificación - are identical
# for every row and stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (values are copied FROM the source
# row's original contents TO the destination row).
$rowMap = @{
    2  = 20
    3  = 6
    4  = 8
    5  = 11
    6  = 31
    7  = 24
    8  = 23
    9  = 25
    10 = 22
    11 = 34
    12 = 26
    13 = 10
    14 = 19
    15 = 17
    16 = 28
    17 = 21
    18 = 30
    19 = 2
    20 = 9
    21 = 5
    22 = 15
    23 = 33
    24 = 27
    25 = 35
    26 = 14
    27 = 16
    28 = 3
    29 = 4
    30 = 29
    31 = 13
    32 = 32
    33 = 12
    34 = 18
    35 = 7
}

# Columns that move together with the observation: D=Fecha, H=Variedad,
# I=Calidad, J=Volumen, K=Precio mínimo, L=Precio máximo,
# M=Precio promedio ponderado, O=Origen, P=Precio $/Kg
$cols = @(4, 8, 9, 10, 11, 12, 13, 15, 16)

# Snapshot the original values for every involved cell before writing
# anything, since rows both read from and are written to.
$snapshot = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $key = "$srcRow-$col"
        if (-not $snapshot.ContainsKey($key)) {
            $snapshot[$key] = $ws.Cells.Item($srcRow, $col).Value2
        }
    }
}

# Now write the permuted values into place.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $key = "$srcRow-$col"
        $ws.Cells.Item($destRow, $col).Value = $snapshot[$key]
    }
}
